$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:D51 to Text format first so numeric-looking values (e.g. "1.010",
# "1.240") keep their exact original string representation (trailing zeros,
# not auto-converted to a number) - matches the source data which is all
# inline/shared text, never numeric cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '28.193.04'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.869.38'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.65%  '
$ws.Range('D5').Value = '313.84'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').Value = '1.009'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('D7').Value = '0.5111'
$ws.Range('E7').Value = '  -0.69%  '
$ws.Range('D8').Value = '0.3891'
$ws.Range('E8').Value = '  +1.09%  '
$ws.Range('D9').Value = '0.08201'
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('D10').Value = '1.113'
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('D11').Value = '41.76'
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').Value = '6.198'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').Value = '1.864.99'
$ws.Range('E13').Value = '  +0.47%  '
$ws.Range('D14').Value = '20.12'
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('D15').Value = '7.183'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').Value = '1.009'
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('D17').Value = '0.00001095'
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').Value = '90.61'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').Value = '0.06659'
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').Value = '17.62'
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('D21').Value = '1.011'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').Value = '5.969'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').Value = '28.246.52'
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('D24').Value = '11.04'
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('D25').Value = '2.246'
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('D26').Value = '2.099.63'
$ws.Range('E26').Value = '  +1.15%  '
$ws.Range('D27').Value = '160.17'
$ws.Range('E27').Value = '  +1.73%  '
$ws.Range('D28').Value = '20.64'
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('D29').Value = '2.401'
$ws.Range('E29').Value = '  -4.54%  '
$ws.Range('D30').Value = '125.57'
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').Value = '0.1043'
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('D32').Value = '1.034'
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').Value = '5.836'
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('D34').Value = '3.616'
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.02419'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').Value = '9.238'
$ws.Range('E36').Value = '  -2.79%  '
$ws.Range('D37').Value = '0.06519'
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('D38').Value = '0.2181'
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '0.6443'
$ws.Range('E39').Value = '  -1.66%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '1.240'
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('D41').Value = '1.174'
$ws.Range('E41').Value = '  -2.45%  '
$ws.Range('D42').Value = '4.935'
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('D43').Value = '11.07'
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('D44').Value = '0.6049'
$ws.Range('E44').Value = '  -1.27%  '
$ws.Range('D45').Value = '13.12'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').Value = '3.688'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').Value = '1.282'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').Value = '1.987'
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('D49').Value = '1.205'
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('D50').Value = '121.61'
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('D51').Value = '0.06893'
$ws.Range('E51').Value = '  +1.05%  '

# Restore the D-column cells to the workbook default ("Normal") style so no
# stray per-cell style index is left behind now that the text is committed.
$priceRange.Style = "Normal"

